# TestData.xlsx / "ScenarioMapping" sheet update.
#
# The SmokeTest column (D) is flipped from "Yes" to "No" for every row
# except TC_009 .. TC_016 (rows 10-17 / the "estateCreation.feature"
# block), which keep "Yes". The RegressionTest column (E) is unchanged
# ("No" throughout).

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("ScenarioMapping")

$ws.Range("D2:D9").Value = "No"
$ws.Range("D18:D335").Value = "No"

# Leave the view roughly where the author left it when saving.
$ws.Range("D338").Select()
